$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.28%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.81%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.098"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08184"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.04%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.110"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.40%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.970"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.40%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.136"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.03%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9287"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.74%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1036"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.27%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1927"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.20%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09170"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.34%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03598"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09902"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001432"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005653"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.62%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.05%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.908"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3414"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1300"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.63%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.102"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.89%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2213"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.27%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04558"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.14%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.63%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004801"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.13%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001252"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.63%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004449"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.27%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01992"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.32%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04943"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.70%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007553"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.04%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1386"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.25%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007906"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.90%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002159"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.67%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01159"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.28%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006602"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.95%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.19%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "200.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "277.27%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-10.43%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.19%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.19%"
